$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 94, pushing existing rows 94..153 down to 95..154
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new price-quote entry
$ws.Range("A94").Value = 7
$ws.Range("B94").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C94").Value = "Ñuble"
$ws.Range("D94").Value = 45001
$ws.Range("E94").Value = 16
$ws.Range("F94").Value = 100112021
$ws.Range("G94").Value = "Ají"
$ws.Range("H94").Value = "Cristal"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 25
$ws.Range("K94").Value = 15000
$ws.Range("L94").Value = 15000
$ws.Range("M94").Value = 15000
$ws.Range("N94").Value = "$/saco 25 kilos"
$ws.Range("O94").Value = "Región del Maule"
$ws.Range("P94").Value = 600
$ws.Range("Q94").Value = 25
$ws.Range("R94").Value = "Hortaliza"
